# Add: Caso2 Completed and Optimized imports
#
# Scenario2 ("Caso2") in the data-driven test workbook is extended from a
# minimal 4-column row (username/password/firstname/lastname) to a full
# completed test-data row covering the whole registration form
# (email, gender, mobile, date-of-birth, subject, hobbies,
#  current-address, state, city, plus the PF first/last name columns),
# and the original placeholder credentials are replaced with the final
# "Test02" data set.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Scenario2")

# New header row (E1:M1) for the additional form fields
$ws2.Range("E1").Value = "email"
$ws2.Range("F1").Value = "gender"
$ws2.Range("G1").Value = "mobile"
$ws2.Range("H1").Value = "date-of-birth"
$ws2.Range("I1").Value = "subject"
$ws2.Range("J1").Value = "hobbies"
$ws2.Range("K1").Value = "current-address"
$ws2.Range("L1").Value = "state"
$ws2.Range("M1").Value = "city"

# Data row (E2:M2) - skip G2/H2 for now, they get special formatting below
$ws2.Range("E2").Value = "prueba@yopmail.com"
$ws2.Range("F2").Value = "Female"

$ws2.Range("I2").Value = "Aspirante"
$ws2.Range("J2").Value = "Music"
$ws2.Range("K2").Value = "Autopista Norte al oriente"
$ws2.Range("L2").Value = "Haryana"
$ws2.Range("M2").Value = "karnal"

# Extra "PF" first/last name columns
$ws2.Range("N2").Value = "Prueba01"
$ws2.Range("O2").Value = "Prueba01"

$ws2.Range("N1").Value = "firstnamePF"
$ws2.Range("O1").Value = "lastnamePF"

# Date of birth is stored as text (numFmtId 49 / "@") so Excel keeps the
# literal dd/mm/yyyy string instead of converting it to a date serial.
$ws2.Range("H2").NumberFormat = "@"
$ws2.Range("H2").Value = "20/09/2000"

# Update the original username/password/firstname/lastname columns to the
# completed "Test02" credential set
$ws2.Range("C2").Value = "Test02"
$ws2.Range("D2").Value = "Test02"

$ws2.Range("A2").Value = "test02"
$ws2.Range("B2").Value = "Test2023*"

# Mobile number stored as a plain number
$ws2.Range("G2").Value = 3004441234
